$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the phone number for the employee in row 4 (I4)
$ws.Range("I4").Value = "+919110790210"

# Update the active selection to J8 (no other data change)
$ws.Range("J8").Select()
